$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle3")

# Row 17
$ws.Range("B17").Value = "without new placement strategy"

# Row 19 (Initial) - set value cell before label so shared-string order matches
$ws.Range("B19").Value = "0=0, 1=48, 2=12, 3=27, 4=124, 5=5, 6=0, 7=0, 8=0, 9=0"
$ws.Range("A19").Value = "Initial"

# Row 18 (Desired)
$ws.Range("B18").Value = "0=0, 1=50, 2=0, 3=13, 4=161, 5=0, 6=0, 7=0, 8=11"
$ws.Range("A18").Value = "Desired"

# Row 20 (Shift)
$ws.Range("B20").Value = "0=0, 1=44, 2=16, 3=27, 4=124, 5=5, 6=0, 7=0, 8=0, 9=0"
$ws.Range("A20").Value = "Shift"

# Row 21 (Break)
$ws.Range("B21").Value = "0=0, 1=43, 2=14, 3=18, 4=136, 5=5, 6=0, 7=0, 8=0, 9=0"
$ws.Range("A21").Value = "Break"

# Row 22 (nUnsatisfied)
$ws.Range("A22").Value = "nUnsatisfied"
$ws.Range("B22").Value = 19

# Row 24 B / Row 23 B
$ws.Range("B24").Value = "{0=0, 1=42, 2=4, 3=24, 4=128, 5=5, 6=0, 7=0, 8=0, 9=0"
$ws.Range("B23").Value = "0=0, 1=59, 2=14, 3=57, 4=68, 5=5, 6=0, 7=0, 8=0, 9=0"

# Row 24 C / Row 23 C
$ws.Range("C24").Value = "0=0, 1=43, 2=14, 3=18, 4=136, 5=5"
$ws.Range("C23").Value = "0=0, 1=48, 2=12, 3=27, 4=124, 5=5"

# Row 24 F / Row 23 F
$ws.Range("F24").Value = "0=0, 1=40, 2=6, 3=24, 4=132, 5=5, 6=0, 7=0, 8=0, 9=0"
$ws.Range("F23").Value = "{0=0, 1=60, 2=26, 3=36, 4=80, 5=5, 6=0, 7=0, 8=0, 9=0}"

# Row 26
$ws.Range("C26").Value = "0=0, 1=40, 2=8, 3=15, 4=136, 5=5"

# Row 27
$ws.Range("C27").Value = "0=0, 1=40, 2=6, 3=21, 4=132, 5=5"

# Update selection to match final state
$ws.Range("F15").Select()
